# [Group18] Danh Gia Thanh Vien Lan 6
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ghi cong")

# Set new scores in column K for rows 16-19
$ws.Range("K16").Value = 1
$ws.Range("K17").Value = 3
$ws.Range("K18").Value = 1
$ws.Range("K19").Value = 3

# Update the active selection on the sheet from I20 to K20
$ws.Range("K20").Select()
